$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows (2-20) of the symbol table were reordered (excludeSymbols were
# filtered out of the list passed into genGigaSymbol, changing the resulting
# row order while keeping each row's data intact). Apply the new row order
# by writing out the full data block in its new arrangement.

$data = @(
    @(601, 9, 60, 67, 60, 42),
    @(801, 3, 67, 65, 52, 45),
    @(1201, 2, 10, 10, 10, 10),
    @(1203, 3, 15, 15, 15, 15),
    @(301, 6, 45, 30, 60, 45),
    @(701, 3, 90, 45, 97, 15),
    @(901, 16, 15, 45, 60, 60),
    @(201, 9, 30, 15, 45, 30),
    @(1202, 2, 10, 10, 10, 10),
    @(101, 9, 30, 15, 60, 15),
    @(1001, 18, 30, 75, 60, 72),
    @(401, 9, 48, 67, 75, 45),
    @(902, 1, 0, 0, 0, 0),
    @(501, 9, 52, 30, 75, 45),
    @(802, 0, 4, 5, 4, 0),
    @(3, 0, 3, 3, 3, 3),
    @(2, 0, 2, 2, 2, 2),
    @(502, 0, 4, 0, 0, 0),
    @(1, 0, 2, 2, 2, 2)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $vals[$col - 1]
    }
}
